# Add columns I (I0) and J (IF) to Sheet1, mirroring the structure of the
# existing H (IP) column: bold/centered/bordered header in row 1 (copied
# from H1's style) and numeric values in rows 2-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style used by the existing header cells (e.g. H1) onto the new
# header cells so they match (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-19
$valuesI = @(8, 9, 8, 1, 7, 9, 5, 9, 8, 8, 9, 4, 6, 7, 5, 9, 7, 8)
$valuesJ = @(8, 9, 8, 1, 7, 9, 5, 9, 8, 8, 9, 5, 6, 7, 6, 9, 7, 8)

for ($i = 0; $i -lt $valuesI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $valuesI[$i]
    $ws.Cells.Item($row, 10).Value = $valuesJ[$i]
}
